$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$cell = $t.Cell(1,1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "20 x 99" + $nl + "  9    9" + $nl + "  ----" + $nl + "2|    |" + $nl + "0|    |"

$cell = $t.Cell(1,2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "97 x 49" + $nl + "  4    9" + $nl + "  ----" + $nl + "9|    |" + $nl + "7|    |"

$cell = $t.Cell(1,3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "16 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "1|    |" + $nl + "6|    |"

$cell = $t.Cell(2,1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "19 x 29" + $nl + "  2    9" + $nl + "  ----" + $nl + "1|    |" + $nl + "9|    |"

$cell = $t.Cell(2,2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "56 x 69" + $nl + "  6    9" + $nl + "  ----" + $nl + "5|    |" + $nl + "6|    |"

$cell = $t.Cell(2,3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "13 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "1|    |" + $nl + "3|    |"

$cell = $t.Cell(3,1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "67 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "6|    |" + $nl + "7|    |"

$cell = $t.Cell(3,2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "58 x 29" + $nl + "  2    9" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"

$cell = $t.Cell(3,3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "63 x 20" + $nl + "  2    0" + $nl + "  ----" + $nl + "6|    |" + $nl + "3|    |"

$cell = $t.Cell(4,1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "78 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "7|    |" + $nl + "8|    |"

$cell = $t.Cell(4,2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "29 x 56" + $nl + "  5    6" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"

$cell = $t.Cell(4,3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "62 x 22" + $nl + "  2    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "2|    |"

$cell = $t.Cell(5,1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "91 x 47" + $nl + "  4    7" + $nl + "  ----" + $nl + "9|    |" + $nl + "1|    |"

$cell = $t.Cell(5,2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "80 x 21" + $nl + "  2    1" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"

$cell = $t.Cell(5,3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "38 x 76" + $nl + "  7    6" + $nl + "  ----" + $nl + "3|    |" + $nl + "8|    |"

